$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper to set a cell's value while keeping it as TEXT (avoids Excel's
# automatic numeric coercion for strings that look like plain numbers,
# e.g. "213.98"). We briefly force a Text number format, assign the
# value, then clear the format again so the cell's style index is left
# exactly as it was before (matches original "no explicit style" cells).
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "25.731.15"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.628.88"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.16%  "

# Row 5 - BNB
Set-TextCell "D5" "213.98"
$ws.Range("E5").Value = "  -0.69%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.82%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.14%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.99%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.55%  "

# Row 10 - Solana
Set-TextCell "D10" "19.44"
$ws.Range("E10").Value = "  -2.06%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.70%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  -0.13%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "1.854.59"
$ws.Range("E13").Value = "  -0.37%  "

# Row 14 - WrappedEther
Set-TextCell "D14" "1.628.63"
$ws.Range("E14").Value = "  -0.58%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.15%  "

# Row 16 - ShibaInu
Set-TextCell "D16" "0.0₃0759"
$ws.Range("E16").Value = "  -2.08%  "

# Row 17 - Litecoin
Set-TextCell "D17" "62.97"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "25.742.34"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.10%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.42%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "191.23"
$ws.Range("E21").Value = "  -1.49%  "

# Row 22 - Avalanche
Set-TextCell "D22" "9.89"
$ws.Range("E22").Value = "  -0.48%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +1.81%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.12%  "

# Row 25 - Toncoin
Set-TextCell "D25" "1.82"
$ws.Range("E25").Value = "  +3.08%  "

# Row 26 - Monero
Set-TextCell "D26" "142.25"
$ws.Range("E26").Value = "  +1.28%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +2.07%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.06%  "

# Row 29 - EthereumClassic
Set-TextCell "D29" "15.44"
$ws.Range("E29").Value = "  -0.74%  "

# Row 30 - PancakeSwap
Set-TextCell "D30" "1.23"
$ws.Range("E30").Value = "  -0.72%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.94%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -0.10%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -1.24%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -1.78%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.70%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  +0.31%  "

# Row 37 - Maker
Set-TextCell "D37" "1.132.97"
$ws.Range("E37").Value = "  +1.68%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  -2.07%  "

# Row 39 - ImmutableX
Set-TextCell "D39" "0.541"
$ws.Range("E39").Value = "  -2.16%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -1.30%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - mCoin
$ws.Range("E42").Value = "  +0.87%  "

# Row 43 - becomes Quant (was FraxShare)
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D43" "100.08"
$ws.Range("E43").Value = "  +1.11%  "

# Row 44 - becomes FraxShare (was Quant)
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D44" "5.51"
$ws.Range("E44").Value = "  -1.13%  "

# Row 45 - TrustWalletToken
$ws.Range("E45").Value = "  -0.62%  "

# Row 46 - RocketPoolETH
Set-TextCell "D46" "1.763.85"
$ws.Range("E46").Value = "  -0.24%  "

# Row 47 - BabyDogeCoin
$ws.Range("E47").Value = "  +0.19%  "

# Row 48 - Aave
Set-TextCell "D48" "55.04"
$ws.Range("E48").Value = "  -0.84%  "

# Row 49 - becomes Mantle (was Cronos)
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D49" "0.418"
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 - becomes Cronos (was Mantle)
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D50" "0.0506"
$ws.Range("E50").Value = "  +0.63%  "

# Row 51 - SynthetixNetwork
$ws.Range("E51").Value = "  -6.94%  "
